# "fix remove odd param" - remove the stray "compliances" column (column F) from
# every order-data sheet. Excel shifts everything after it one column to the
# left automatically; the header cell (shared string "compliances" /
# "Edit compliances" on FPA004-006-010) disappears from the sheet entirely.

$wb = $excel.ActiveWorkbook

# 1. FPA001
$ws = $wb.Worksheets.Item("FPA001")
$ws.Columns("F:F").Delete()
$ws.Range("F24").Select()

# 2. FPA002-003-005-007
$ws = $wb.Worksheets.Item("FPA002-003-005-007")
$ws.Columns("F:F").Delete()
$ws.Range("F1:F1048576").Select()

# 3. FPA004-006-010
$ws = $wb.Worksheets.Item("FPA004-006-010")
$ws.Columns("F:F").Delete()
$ws.Range("F29").Select()

# 4. FPA008-009
$ws = $wb.Worksheets.Item("FPA008-009")
$ws.Columns("F:F").Delete()
$ws.Range("G27").Select()
# The deleted column sat inside the hidden _FilterDatabase defined name for
# this sheet; Excel shrinks that range along with the column delete.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlnm._FilterDatabase" -or $n.Name -like "*FilterDatabase*") {
        $n.RefersTo = "='FPA008-009'!`$A`$1:`$O`$4"
    }
}

# 5. BTMI002
$ws = $wb.Worksheets.Item("BTMI002")
$ws.Columns("F:F").Delete()
$ws.Range("H45").Select()

# 6. BTMI003
$ws = $wb.Worksheets.Item("BTMI003")
$ws.Columns("F:F").Delete()
$ws.Range("H24").Select()

# 7. BTMI015 - ends up the active sheet/tab
$ws = $wb.Worksheets.Item("BTMI015")
$ws.Columns("F:F").Delete()
$ws.Activate()
$ws.Range("J32").Select()
